# Regenerate save_data column G ("K") values: use K instead of Strike#,
# regen std/mean, calc and write s_vals.
#
# This updates the "K" column (column G) for rows 2-20 on the active sheet
# with the newly computed strike-count values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newKValues = [ordered]@{
    2  = 2
    3  = 2
    4  = 1
    5  = 6
    6  = 2
    7  = 1
    8  = 3
    9  = 3
    10 = 3
    11 = 3
    12 = 0
    13 = 4
    14 = 0
    15 = 5
    16 = 7
    17 = 6
    18 = 4
    19 = 1
    20 = 1
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
